$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 96 - this shifts the old rows 96-100 down to 97-101
$ws.Rows.Item(96).Insert()

# Populate the newly inserted row 96 with the new weekly data point
$ws.Cells.Item(96, 1).Value = 1
$ws.Cells.Item(96, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(96, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(96, 4).Value = 44610
$ws.Cells.Item(96, 5).Value = 15
$ws.Cells.Item(96, 6).Value = 100112042
$ws.Cells.Item(96, 7).Value = "Locoto"
$ws.Cells.Item(96, 8).Value = "Sin especificar"
$ws.Cells.Item(96, 9).Value = "Primera"
$ws.Cells.Item(96, 10).Value = 120
$ws.Cells.Item(96, 11).Value = 37000
$ws.Cells.Item(96, 12).Value = 38000
$ws.Cells.Item(96, 13).Value = 37500
$ws.Cells.Item(96, 14).Value = "$/caja 20 kilos"
$ws.Cells.Item(96, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(96, 16).Value = 1875
$ws.Cells.Item(96, 17).Value = 20
$ws.Cells.Item(96, 18).Value = "Hortaliza"
